$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 568.7
$ws.Range("I98").Value = 581.5263
$ws.Range("J98").Value = 325
$ws.Range("K98").Value = 581.5263
$ws.Range("L98").Value = 325
$ws.Range("M98").Value = 916.4737
$ws.Range("N98").Value = -3321
$ws.Range("H113").Value = 6701.9414
$ws.Range("I113").Value = 2978
$ws.Range("K113").Value = 2978
$ws.Range("M113").Value = 276
$ws.Range("H122").Value = 568.7
$ws.Range("I122").Value = 581.5263
$ws.Range("J122").Value = 325
$ws.Range("K122").Value = 1744.5789
$ws.Range("L122").Value = 975
$ws.Range("M122").Value = 705.4211
$ws.Range("N122").Value = -5875
$ws.Range("H129").Value = 3205.182
$ws.Range("I129").Value = 3769
$ws.Range("K129").Value = 11307
$ws.Range("M129").Value = -6307
$ws.Range("H137").Value = 3198.6667
$ws.Range("I137").Value = 2798
$ws.Range("J137").Value = 3399
$ws.Range("K137").Value = 8394
$ws.Range("L137").Value = 10197
$ws.Range("M137").Value = -5844
$ws.Range("N137").Value = -15297

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1972.5
$ws.Range("I2").Value = 1296.6666
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 1296.6666
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -1183.6666
$ws.Range("N2").Value = -4226
$ws.Range("H32").Value = 3152.4443
$ws.Range("I32").Value = 1807.4375
$ws.Range("K32").Value = 1807.4375
$ws.Range("M32").Value = -1520.4375
$ws.Range("H45").Value = 2011.7262
$ws.Range("I45").Value = 2011.7262
$ws.Range("K45").Value = 2011.7262
$ws.Range("M45").Value = -1634.7262
$ws.Range("H61").Value = 4317.4634
$ws.Range("I61").Value = 4104.8965
$ws.Range("J61").Value = 4831.1665
$ws.Range("K61").Value = 4104.8965
$ws.Range("L61").Value = 4831.1665
$ws.Range("M61").Value = -3892.8965
$ws.Range("N61").Value = -5255.1665
$ws.Range("H74").Value = 2731.6667
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 2731.6667
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 2731.6667
$ws.Range("N74").Value = -4479.6667
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 2731.6667
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 2731.6667
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 13658.3335
$ws.Range("N77").Value = -22394.3335
$ws.Range("M77").ClearContents()
$ws.Range("H102").Value = 6949.6665
$ws.Range("I102").Value = 5899.3335
$ws.Range("K102").Value = 5899.3335
$ws.Range("M102").Value = -4277.3335
$ws.Range("H116").Value = 1972.5
$ws.Range("I116").Value = 1296.6666
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 1296.6666
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 997.3334
$ws.Range("N116").Value = -8588
$ws.Range("H122").Value = 3748.625
$ws.Range("I122").Value = 3712.8572
$ws.Range("K122").Value = 11138.5716
$ws.Range("M122").Value = -8688.571599999999
$ws.Range("H136").Value = 4317.4634
$ws.Range("I136").Value = 4104.8965
$ws.Range("J136").Value = 4831.1665
$ws.Range("K136").Value = 12314.6895
$ws.Range("L136").Value = 14493.4995
$ws.Range("M136").Value = -9764.6895
$ws.Range("N136").Value = -19593.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1972.5
$ws.Range("I3").Value = 1296.6666
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 1296.6666
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -1182.6666
$ws.Range("N3").Value = -4228
$ws.Range("H86").Value = 1744.909
$ws.Range("J86").Value = 970.3333
$ws.Range("L86").Value = 970.3333
$ws.Range("N86").Value = -3216.3333
$ws.Range("H89").Value = 1744.909
$ws.Range("J89").Value = 970.3333
$ws.Range("L89").Value = 4851.6665
$ws.Range("N89").Value = -16083.6665
$ws.Range("H134").Value = 2012.4517
$ws.Range("I134").Value = 2311.9412
$ws.Range("J134").Value = 1648.7858
$ws.Range("K134").Value = 6935.823600000001
$ws.Range("L134").Value = 4946.357400000001
$ws.Range("M134").Value = -4400.823600000001
$ws.Range("N134").Value = -10016.3574
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 665.6667
$ws.Range("I16").Value = 665.6667
$ws.Range("K16").Value = 665.6667
$ws.Range("M16").Value = -378.6667
$ws.Range("H31").Value = 3706.25
$ws.Range("J31").Value = 4382
$ws.Range("L31").Value = 4382
$ws.Range("N31").Value = -4972
$ws.Range("H34").Value = 3706.25
$ws.Range("J34").Value = 4382
$ws.Range("L34").Value = 4382
$ws.Range("N34").Value = -4786
$ws.Range("H86").Value = 6118.8335
$ws.Range("J86").Value = 5676.5
$ws.Range("L86").Value = 5676.5
$ws.Range("N86").Value = -7922.5
$ws.Range("H89").Value = 6118.8335
$ws.Range("J89").Value = 5676.5
$ws.Range("L89").Value = 28382.5
$ws.Range("N89").Value = -39614.5
$ws.Range("H99").Value = 17712.615
$ws.Range("I99").Value = 13229.444
$ws.Range("K99").Value = 13229.444
$ws.Range("M99").Value = -11731.444
$ws.Range("H107").Value = 395.85715
$ws.Range("I107").Value = 380.80487
$ws.Range("K107").Value = 380.80487
$ws.Range("M107").Value = 1539.19513
$ws.Range("H113").Value = 665.6667
$ws.Range("I113").Value = 665.6667
$ws.Range("K113").Value = 665.6667
$ws.Range("M113").Value = 1504.3333
$ws.Range("H126").Value = 17712.615
$ws.Range("I126").Value = 13229.444
$ws.Range("K126").Value = 39688.33199999999
$ws.Range("M126").Value = -37218.33199999999
$ws.Range("H132").Value = 4364.4287
$ws.Range("I132").Value = 3799.6
$ws.Range("K132").Value = 11398.8
$ws.Range("M132").Value = -8868.799999999999
$ws.Range("H134").Value = 4869.6
$ws.Range("I134").Value = 4072.25
$ws.Range("J134").Value = 5780.857
$ws.Range("K134").Value = 12216.75
$ws.Range("L134").Value = 17342.571
$ws.Range("M134").Value = -9681.75
$ws.Range("N134").Value = -22412.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 19872
$ws.Range("J92").Value = 19872
$ws.Range("L92").Value = 19872
$ws.Range("N92").Value = -23616
$ws.Range("H102").Value = 2676.4167
$ws.Range("I102").Value = 2837.5557
$ws.Range("K102").Value = 2837.5557
$ws.Range("M102").Value = -1215.5557
$ws.Range("H107").Value = 1684.5555
$ws.Range("I107").Value = 744
$ws.Range("J107").Value = 2154.8333
$ws.Range("K107").Value = 744
$ws.Range("L107").Value = 2154.8333
$ws.Range("M107").Value = 1176
$ws.Range("N107").Value = -5994.8333
$ws.Range("H122").Value = 1701.875
$ws.Range("I122").Value = 1590.7142
$ws.Range("K122").Value = 4772.142599999999
$ws.Range("M122").Value = -2322.142599999999
$ws.Range("H132").Value = 5337
$ws.Range("I132").Value = 3011
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 9033
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -6503
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7198.722
$ws.Range("I61").Value = 6995.364
$ws.Range("J61").Value = 7518.2856
$ws.Range("K61").Value = 6995.364
$ws.Range("L61").Value = 7518.2856
$ws.Range("M61").Value = -6793.364
$ws.Range("N61").Value = -7922.2856
$ws.Range("H100").Value = 3713.8572
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 7198.722
$ws.Range("I113").Value = 6995.364
$ws.Range("J113").Value = 7518.2856
$ws.Range("K113").Value = 6995.364
$ws.Range("L113").Value = 7518.2856
$ws.Range("M113").Value = -4825.364
$ws.Range("N113").Value = -11858.2856
$ws.Range("H132").Value = 2927.4443
$ws.Range("I132").Value = 3024.625
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 9073.875
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -6543.875
$ws.Range("N132").Value = -11510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4194
$ws.Range("I62").Value = 4194
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4194
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3570
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4194
$ws.Range("I65").Value = 4194
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20970
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -17850
$ws.Range("N65").ClearContents()
$ws.Range("H126").Value = 4810.2
$ws.Range("I126").Value = 3229.6667
$ws.Range("K126").Value = 9689.000100000001
$ws.Range("M126").Value = -7219.000100000001
$ws.Range("H132").Value = 8231.207
$ws.Range("I132").Value = 5874.478
$ws.Range("J132").Value = 17265.334
$ws.Range("K132").Value = 17623.434
$ws.Range("L132").Value = 51796.00199999999
$ws.Range("M132").Value = -15093.434
$ws.Range("N132").Value = -56856.00199999999
$ws.Range("H136").Value = 2441.0715
$ws.Range("I136").Value = 1178.1
$ws.Range("J136").Value = 5598.5
$ws.Range("K136").Value = 3534.3
$ws.Range("L136").Value = 16795.5
$ws.Range("M136").Value = -984.2999999999997
$ws.Range("N136").Value = -21895.5
